$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at D (shifts existing D:K data to F:M)
$ws.Range("D:E").Insert()

# Fix formatting of the newly inserted D:E columns.
# Default data-row format (style matches column F, e.g. F8)
$ws.Range("F8").Copy()
$ws.Range("D7:E102").PasteSpecial(-4122)

# Date-row format (rows 7, 38, 80 use a date number format, matching column F on those rows)
$ws.Range("F7").Copy()
$ws.Range("D7:E7").PasteSpecial(-4122)
$ws.Range("F38").Copy()
$ws.Range("D38:E38").PasteSpecial(-4122)
$ws.Range("F80").Copy()
$ws.Range("D80:E80").PasteSpecial(-4122)

# Rows 36 and 78 are blank spacer rows with no cells in either the source or
# target layout; undo the stray format-only cells the bulk paste above left there.
$ws.Range("D36:E36").Clear()
$ws.Range("D78:E78").Clear()

$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("F7").Value = 43281
$ws.Range("G7").Value = 43190
$ws.Range("H7").Value = 43100
$ws.Range("I7").Value = 43008
$ws.Range("J7").Value = 42916
$ws.Range("K7").Value = 42825
$ws.Range("L7").Value = 42735
$ws.Range("M7").Value = 42643
$ws.Range("D8").Value = 1074800
$ws.Range("E8").Value = 909500
$ws.Range("F8").Value = 940800
$ws.Range("G8").Value = 1183100
$ws.Range("H8").Value = 946100
$ws.Range("I8").Value = 922600
$ws.Range("J8").Value = 810200
$ws.Range("K8").Value = 969300
$ws.Range("L8").Value = 924100
$ws.Range("M8").Value = 766600
$ws.Range("D9").Value = 677100
$ws.Range("E9").Value = 502700
$ws.Range("F9").Value = 550100
$ws.Range("G9").Value = 781900
$ws.Range("H9").Value = 532500
$ws.Range("I9").Value = 428100
$ws.Range("J9").Value = 389300
$ws.Range("K9").Value = 575000
$ws.Range("L9").Value = 487500
$ws.Range("M9").Value = 320500
$ws.Range("D10").Value = 397700
$ws.Range("E10").Value = 406800
$ws.Range("F10").Value = 390700
$ws.Range("G10").Value = 401200
$ws.Range("H10").Value = 413600
$ws.Range("I10").Value = 494500
$ws.Range("J10").Value = 420900
$ws.Range("K10").Value = 394300
$ws.Range("L10").Value = 436600
$ws.Range("M10").Value = 446100
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("F12").Value = "NA"
$ws.Range("G12").Value = "NA"
$ws.Range("H12").Value = "NA"
$ws.Range("I12").Value = "NA"
$ws.Range("J12").Value = "NA"
$ws.Range("K12").Value = "NA"
$ws.Range("L12").Value = "NA"
$ws.Range("M12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 537000
$ws.Range("F14").Value = "NA"
$ws.Range("G14").Value = "NA"
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = "NA"
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = 0
$ws.Range("D15").Value = 70300
$ws.Range("E15").Value = 68500
$ws.Range("F15").Value = 66600
$ws.Range("G15").Value = 64100
$ws.Range("H15").Value = 73300
$ws.Range("I15").Value = 65700
$ws.Range("J15").Value = 64800
$ws.Range("K15").Value = 65500
$ws.Range("L15").Value = 66400
$ws.Range("M15").Value = 63500
$ws.Range("D17").Value = 581100
$ws.Range("E17").Value = 1301400
$ws.Range("F17").Value = 795100
$ws.Range("G17").Value = 1013600
$ws.Range("H17").Value = 772300
$ws.Range("I17").Value = 754300
$ws.Range("J17").Value = 639400
$ws.Range("K17").Value = 796700
$ws.Range("L17").Value = 767100
$ws.Range("M17").Value = 560400
$ws.Range("D18").Value = 493700
$ws.Range("E18").Value = -391900
$ws.Range("F18").Value = 145700
$ws.Range("G18").Value = 169500
$ws.Range("H18").Value = 173800
$ws.Range("I18").Value = 168300
$ws.Range("J18").Value = 170800
$ws.Range("K18").Value = 172600
$ws.Range("L18").Value = 157000
$ws.Range("M18").Value = 206200
$ws.Range("D20").Value = 2300
$ws.Range("E20").Value = -292500
$ws.Range("F20").Value = 8000
$ws.Range("G20").Value = 7200
$ws.Range("H20").Value = 12900
$ws.Range("I20").Value = 8900
$ws.Range("J20").Value = 3100
$ws.Range("K20").Value = 9800
$ws.Range("L20").Value = 3300
$ws.Range("M20").Value = 2800
$ws.Range("D21").Value = 566300
$ws.Range("E21").Value = -615900
$ws.Range("F21").Value = 220200
$ws.Range("G21").Value = 240800
$ws.Range("H21").Value = 259900
$ws.Range("I21").Value = 242900
$ws.Range("J21").Value = 238700
$ws.Range("K21").Value = 247900
$ws.Range("L21").Value = 226800
$ws.Range("M21").Value = 272500
$ws.Range("D22").Value = 11200
$ws.Range("E22").Value = 60300
$ws.Range("F22").Value = 59600
$ws.Range("G22").Value = 59100
$ws.Range("H22").Value = 56700
$ws.Range("I22").Value = 56600
$ws.Range("J22").Value = 56400
$ws.Range("K22").Value = 55900
$ws.Range("L22").Value = 50800
$ws.Range("M22").Value = 48500
$ws.Range("D23").Value = 484800
$ws.Range("E23").Value = -744700
$ws.Range("F23").Value = 94100
$ws.Range("G23").Value = 117600
$ws.Range("H23").Value = 129900
$ws.Range("I23").Value = 120700
$ws.Range("J23").Value = 117400
$ws.Range("K23").Value = 126500
$ws.Range("L23").Value = 109500
$ws.Range("M23").Value = 160600
$ws.Range("D24").Value = 1800
$ws.Range("E24").Value = 600
$ws.Range("F24").Value = 800
$ws.Range("G24").Value = 500
$ws.Range("H24").Value = -800
$ws.Range("I24").Value = 400
$ws.Range("J24").Value = 1000
$ws.Range("K24").Value = 200
$ws.Range("L24").Value = 600
$ws.Range("M24").Value = 300
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = 0
$ws.Range("D26").Value = 483000
$ws.Range("E26").Value = -745300
$ws.Range("F26").Value = 93300
$ws.Range("G26").Value = 117100
$ws.Range("H26").Value = 130800
$ws.Range("I26").Value = 120200
$ws.Range("J26").Value = 116400
$ws.Range("K26").Value = 126300
$ws.Range("L26").Value = 108900
$ws.Range("M26").Value = 160300
$ws.Range("D27").Value = 482500
$ws.Range("E27").Value = -745800
$ws.Range("F27").Value = 91900
$ws.Range("G27").Value = 112400
$ws.Range("H27").Value = 126300
$ws.Range("I27").Value = 116200
$ws.Range("J27").Value = 112700
$ws.Range("K27").Value = 123600
$ws.Range("L27").Value = 107700
$ws.Range("M27").Value = 156400
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 0
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = 0
$ws.Range("G29").Value = 0
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 0
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = 0
$ws.Range("D32").Value = -2300
$ws.Range("E32").Value = 292500
$ws.Range("F32").Value = -8000
$ws.Range("G32").Value = -7200
$ws.Range("H32").Value = -12900
$ws.Range("I32").Value = -8900
$ws.Range("J32").Value = -3100
$ws.Range("K32").Value = -9800
$ws.Range("L32").Value = -3300
$ws.Range("M32").Value = -2800
$ws.Range("D33").Value = 482500
$ws.Range("E33").Value = -745800
$ws.Range("F33").Value = 91900
$ws.Range("G33").Value = 112400
$ws.Range("H33").Value = 126300
$ws.Range("I33").Value = 116200
$ws.Range("J33").Value = 112700
$ws.Range("K33").Value = 123600
$ws.Range("L33").Value = 107700
$ws.Range("M33").Value = 156400
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 0
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = 0
$ws.Range("D35").Value = 482500
$ws.Range("E35").Value = -745800
$ws.Range("F35").Value = 91900
$ws.Range("G35").Value = 112400
$ws.Range("H35").Value = 126300
$ws.Range("I35").Value = 116200
$ws.Range("J35").Value = 112700
$ws.Range("K35").Value = 123600
$ws.Range("L35").Value = 107700
$ws.Range("M35").Value = 156400
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("F38").Value = 43281
$ws.Range("G38").Value = 43190
$ws.Range("H38").Value = 43100
$ws.Range("I38").Value = 43008
$ws.Range("J38").Value = 42916
$ws.Range("K38").Value = 42825
$ws.Range("L38").Value = 42735
$ws.Range("M38").Value = 42643
$ws.Range("D41").Value = 1800
$ws.Range("E41").Value = 700
$ws.Range("F41").Value = 1700
$ws.Range("G41").Value = 6600
$ws.Range("H41").Value = 2200
$ws.Range("I41").Value = 7900
$ws.Range("J41").Value = 2400
$ws.Range("K41").Value = 4900
$ws.Range("L41").Value = 640300
$ws.Range("M41").Value = 0
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("F42").Value = 0
$ws.Range("G42").Value = 0
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = 0
$ws.Range("D43").Value = 251000
$ws.Range("E43").Value = 278500
$ws.Range("F43").Value = 256100
$ws.Range("G43").Value = 287300
$ws.Range("H43").Value = 295700
$ws.Range("I43").Value = 234600
$ws.Range("J43").Value = 219800
$ws.Range("K43").Value = 229800
$ws.Range("L43").Value = 255300
$ws.Range("M43").Value = 273200
$ws.Range("D44").Value = 210900
$ws.Range("E44").Value = 197500
$ws.Range("F44").Value = 155900
$ws.Range("G44").Value = 158200
$ws.Range("H44").Value = 301400
$ws.Range("I44").Value = 233400
$ws.Range("J44").Value = 286500
$ws.Range("K44").Value = 305100
$ws.Range("L44").Value = 356800
$ws.Range("M44").Value = 295300
$ws.Range("D45").Value = 59900
$ws.Range("E45").Value = 101600
$ws.Range("F45").Value = 93300
$ws.Range("G45").Value = 67900
$ws.Range("H45").Value = 57300
$ws.Range("I45").Value = 82400
$ws.Range("J45").Value = 66900
$ws.Range("K45").Value = 65600
$ws.Range("L45").Value = 66400
$ws.Range("M45").Value = 72400
$ws.Range("D46").Value = 523600
$ws.Range("E46").Value = 578200
$ws.Range("F46").Value = 507000
$ws.Range("G46").Value = 520000
$ws.Range("H46").Value = 656600
$ws.Range("I46").Value = 558300
$ws.Range("J46").Value = 575500
$ws.Range("K46").Value = 605400
$ws.Range("L46").Value = 1318900
$ws.Range("M46").Value = 641000
$ws.Range("D47").Value = 1113900
$ws.Range("E47").Value = 1168300
$ws.Range("F47").Value = 1482100
$ws.Range("G47").Value = 1496600
$ws.Range("H47").Value = 1494400
$ws.Range("I47").Value = 1496000
$ws.Range("J47").Value = 1254800
$ws.Range("K47").Value = 1249600
$ws.Range("L47").Value = 89600
$ws.Range("M47").Value = 89900
$ws.Range("D48").Value = 6914800
$ws.Range("E48").Value = 6973400
$ws.Range("F48").Value = 6889400
$ws.Range("G48").Value = 6806200
$ws.Range("H48").Value = 6735800
$ws.Range("I48").Value = 6633200
$ws.Range("J48").Value = 6579100
$ws.Range("K48").Value = 6524200
$ws.Range("L48").Value = 6483300
$ws.Range("M48").Value = 6422000
$ws.Range("D49").Value = 764200
$ws.Range("E49").Value = 780300
$ws.Range("F49").Value = 1333800
$ws.Range("G49").Value = 1350200
$ws.Range("H49").Value = 1366400
$ws.Range("I49").Value = 1378200
$ws.Range("J49").Value = 1394900
$ws.Range("K49").Value = 1411500
$ws.Range("L49").Value = 1427800
$ws.Range("M49").Value = 1441500
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("F50").Value = 0
$ws.Range("G50").Value = 0
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("M50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = 0
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = 0
$ws.Range("D52").Value = 39100
$ws.Range("E52").Value = 35700
$ws.Range("F52").Value = 38600
$ws.Range("G52").Value = 42500
$ws.Range("H52").Value = 51500
$ws.Range("I52").Value = 46800
$ws.Range("J52").Value = 76400
$ws.Range("K52").Value = 76200
$ws.Range("L52").Value = 101500
$ws.Range("M52").Value = 42700
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("M53").Value = 0
$ws.Range("D54").Value = 9355600
$ws.Range("E54").Value = 9535900
$ws.Range("F54").Value = 10250900
$ws.Range("G54").Value = 10215500
$ws.Range("H54").Value = 10304700
$ws.Range("I54").Value = 10112500
$ws.Range("J54").Value = 9880700
$ws.Range("K54").Value = 9866800
$ws.Range("L54").Value = 9421100
$ws.Range("M54").Value = 8637000
$ws.Range("D57").Value = 147500
$ws.Range("E57").Value = 145000
$ws.Range("F57").Value = 86300
$ws.Range("G57").Value = 92900
$ws.Range("H57").Value = 160800
$ws.Range("I57").Value = 82400
$ws.Range("J57").Value = 64900
$ws.Range("K57").Value = 55300
$ws.Range("L57").Value = 107400
$ws.Range("M57").Value = 91900
$ws.Range("D58").Value = 702700
$ws.Range("E58").Value = 178200
$ws.Range("F58").Value = 171000
$ws.Range("G58").Value = 115200
$ws.Range("H58").Value = 252200
$ws.Range("I58").Value = 185400
$ws.Range("J58").Value = 227300
$ws.Range("K58").Value = 258300
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = 218300
$ws.Range("D59").Value = 277500
$ws.Range("E59").Value = 279100
$ws.Range("F59").Value = 245600
$ws.Range("G59").Value = 238100
$ws.Range("H59").Value = 272400
$ws.Range("I59").Value = 241500
$ws.Range("J59").Value = 246200
$ws.Range("K59").Value = 234200
$ws.Range("L59").Value = 292200
$ws.Range("M59").Value = 275200
$ws.Range("D60").Value = 1127800
$ws.Range("E60").Value = 602400
$ws.Range("F60").Value = 502900
$ws.Range("G60").Value = 446200
$ws.Range("H60").Value = 685400
$ws.Range("I60").Value = 509300
$ws.Range("J60").Value = 538400
$ws.Range("K60").Value = 547800
$ws.Range("L60").Value = 399500
$ws.Range("M60").Value = 585400
$ws.Range("D61").Value = 4011700
$ws.Range("E61").Value = 4985200
$ws.Range("F61").Value = 4877900
$ws.Range("G61").Value = 4587900
$ws.Range("H61").Value = 4658300
$ws.Range("I61").Value = 4593600
$ws.Range("J61").Value = 4579900
$ws.Range("K61").Value = 4561000
$ws.Range("L61").Value = 4217700
$ws.Range("M61").Value = 3826900
$ws.Range("D62").Value = 84900
$ws.Range("E62").Value = 88200
$ws.Range("F62").Value = 90500
$ws.Range("G62").Value = 92200
$ws.Range("H62").Value = 92700
$ws.Range("I62").Value = 96600
$ws.Range("J62").Value = 98900
$ws.Range("K62").Value = 98500
$ws.Range("L62").Value = 105400
$ws.Range("M62").Value = 111800
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("F63").Value = 0
$ws.Range("G63").Value = 0
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("F64").Value = 0
$ws.Range("G64").Value = 0
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("F65").Value = 0
$ws.Range("G65").Value = 0
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = 0
$ws.Range("D66").Value = 5235500
$ws.Range("E66").Value = 5687000
$ws.Range("F66").Value = 5482600
$ws.Range("G66").Value = 5405800
$ws.Range("H66").Value = 5713700
$ws.Range("I66").Value = 5482000
$ws.Range("J66").Value = 5500500
$ws.Range("K66").Value = 5494300
$ws.Range("L66").Value = 5009400
$ws.Range("M66").Value = 4811800
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("F68").Value = 0
$ws.Range("G68").Value = 0
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("F69").Value = 0
$ws.Range("G69").Value = 0
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("F70").Value = 0
$ws.Range("G70").Value = 0
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("F71").Value = 0
$ws.Range("G71").Value = 0
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = 0
$ws.Range("D72").Value = 0
$ws.Range("E72").Value = 0
$ws.Range("F72").Value = 0
$ws.Range("G72").Value = 0
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = 0
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("F73").Value = 0
$ws.Range("G73").Value = 0
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("F74").Value = 0
$ws.Range("G74").Value = 0
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("F75").Value = 0
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = 0
$ws.Range("D76").Value = 4120100
$ws.Range("E76").Value = 3848900
$ws.Range("F76").Value = 4768300
$ws.Range("G76").Value = 4809700
$ws.Range("H76").Value = 4590900
$ws.Range("I76").Value = 4630500
$ws.Range("J76").Value = 4380200
$ws.Range("K76").Value = 4372500
$ws.Range("L76").Value = 4411700
$ws.Range("M76").Value = 3825200
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("F77").Value = 0
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("F80").Value = 43281
$ws.Range("G80").Value = 43190
$ws.Range("H80").Value = 43100
$ws.Range("I80").Value = 43008
$ws.Range("J80").Value = 42916
$ws.Range("K80").Value = 42825
$ws.Range("L80").Value = 42735
$ws.Range("M80").Value = 42643
$ws.Range("D81").Value = 482500
$ws.Range("E81").Value = -745800
$ws.Range("F81").Value = 91900
$ws.Range("G81").Value = 112400
$ws.Range("H81").Value = 126300
$ws.Range("I81").Value = 116200
$ws.Range("J81").Value = 112700
$ws.Range("K81").Value = 123600
$ws.Range("L81").Value = 107700
$ws.Range("M81").Value = 156400
$ws.Range("D83").Value = 70300
$ws.Range("E83").Value = 68500
$ws.Range("F83").Value = 66600
$ws.Range("G83").Value = 64100
$ws.Range("H83").Value = 73300
$ws.Range("I83").Value = 65700
$ws.Range("J83").Value = 64800
$ws.Range("K83").Value = 65500
$ws.Range("L83").Value = 66400
$ws.Range("M83").Value = 63500
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("F84").Value = 0
$ws.Range("G84").Value = 0
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("F85").Value = 0
$ws.Range("G85").Value = 0
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("F86").Value = 0
$ws.Range("G86").Value = 0
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("F87").Value = 0
$ws.Range("G87").Value = 0
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("M87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("F88").Value = 0
$ws.Range("G88").Value = 0
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = 0
$ws.Range("D89").Value = 171600
$ws.Range("E89").Value = 165000
$ws.Range("F89").Value = 195000
$ws.Range("G89").Value = 236200
$ws.Range("H89").Value = 201100
$ws.Range("I89").Value = 251100
$ws.Range("J89").Value = 243200
$ws.Range("K89").Value = 193000
$ws.Range("L89").Value = 168100
$ws.Range("M89").Value = 178200
$ws.Range("D91").Value = -110600
$ws.Range("E91").Value = -107000
$ws.Range("F91").Value = -133100
$ws.Range("G91").Value = -116900
$ws.Range("H91").Value = -129600
$ws.Range("I91").Value = -100700
$ws.Range("J91").Value = -104700
$ws.Range("K91").Value = -98300
$ws.Range("L91").Value = -106200
$ws.Range("M91").Value = -156000
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("F92").Value = 0
$ws.Range("G92").Value = 0
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("F93").Value = 0
$ws.Range("G93").Value = 0
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = 0
$ws.Range("D94").Value = 397900
$ws.Range("E94").Value = -90300
$ws.Range("F94").Value = -146100
$ws.Range("G94").Value = -89100
$ws.Range("H94").Value = -140500
$ws.Range("I94").Value = -329000
$ws.Range("J94").Value = -93300
$ws.Range("K94").Value = -1247200
$ws.Range("L94").Value = -103500
$ws.Range("M94").Value = -175700
$ws.Range("D96").Value = -115600
$ws.Range("E96").Value = -186000
$ws.Range("F96").Value = -186100
$ws.Range("G96").Value = -186100
$ws.Range("H96").Value = -185300
$ws.Range("I96").Value = -178900
$ws.Range("J96").Value = -176100
$ws.Range("K96").Value = -174100
$ws.Range("L96").Value = -171900
$ws.Range("M96").Value = -159100
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("F97").Value = 0
$ws.Range("G97").Value = 0
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("F98").Value = 0
$ws.Range("G98").Value = 0
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("F99").Value = 0
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 0
$ws.Range("D100").Value = -568400
$ws.Range("E100").Value = -75800
$ws.Range("F100").Value = -53800
$ws.Range("G100").Value = -142700
$ws.Range("H100").Value = -66400
$ws.Range("I100").Value = 83500
$ws.Range("J100").Value = -152300
$ws.Range("K100").Value = 418700
$ws.Range("L100").Value = 575700
$ws.Range("M100").Value = -17200
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("F101").Value = 0
$ws.Range("G101").Value = 0
$ws.Range("H101").Value = 0
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("M101").Value = 0
$ws.Range("D102").Value = 1200
$ws.Range("E102").Value = -1000
$ws.Range("F102").Value = -5000
$ws.Range("G102").Value = 4400
$ws.Range("H102").Value = -5700
$ws.Range("I102").Value = 5500
$ws.Range("J102").Value = -2500
$ws.Range("K102").Value = -635500
$ws.Range("L102").Value = 640300
$ws.Range("M102").Value = -14800
